## Fruta / hortaliza, semanal
## Insert a new weekly price-report row for "Ají" (Hortaliza) at row 122,
## pushing the existing rows 122-211 down to 123-212.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 122. This shifts every
# row from 122 through 211 down by one (to 123-212), keeping their data,
# styles and formatting intact - exactly matching the diff, which shows
# each row N (122..211) taking on the previous values of row N-1, and the
# former row 211 becoming the new row 212.
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with the new data point.
$ws.Range("A122").Value() = 7
$ws.Range("B122").Value() = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C122").Value() = "Ñuble"
$ws.Range("D122").Value() = 45072
$ws.Range("E122").Value() = 16
$ws.Range("F122").Value() = 100112021
$ws.Range("G122").Value() = "Ají"
$ws.Range("H122").Value() = "Cacho cabra rojo"
$ws.Range("I122").Value() = "Primera"
$ws.Range("J122").Value() = 40
$ws.Range("K122").Value() = 16000
$ws.Range("L122").Value() = 17000
$ws.Range("M122").Value() = 16500
$ws.Range("N122").Value() = "`$/saco 25 kilos"
$ws.Range("O122").Value() = "Región del Maule"
$ws.Range("P122").Value() = 660
$ws.Range("Q122").Value() = 25
$ws.Range("R122").Value() = "Hortaliza"
